# Insert a new weekly price record for "Femacal de La Calera - Zanahoria"
# at row 515. This shifts the existing rows 515:540 down to 516:541 and
# keeps all of their data intact (Excel's native row-insert semantics).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(515).Insert()

$ws.Cells.Item(515, 1).Value = 3
$ws.Cells.Item(515, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(515, 3).Value = "Coquimbo"
$ws.Cells.Item(515, 4).Value = 45041
$ws.Cells.Item(515, 5).Value = 5
$ws.Cells.Item(515, 6).Value = 100114013
$ws.Cells.Item(515, 7).Value = "Zanahoria"
$ws.Cells.Item(515, 8).Value = "Sin especificar"
$ws.Cells.Item(515, 9).Value = "Primera"
$ws.Cells.Item(515, 10).Value = 310
$ws.Cells.Item(515, 11).Value = 7000
$ws.Cells.Item(515, 12).Value = 8000
$ws.Cells.Item(515, 13).Value = 7484
$ws.Cells.Item(515, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(515, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(515, 16).Value = 374
$ws.Cells.Item(515, 17).Value = 20
$ws.Cells.Item(515, 18).Value = "Hortaliza"

# Column D carries the date number format used throughout the sheet;
# make sure the newly inserted row's date cell keeps that same format.
$ws.Cells.Item(515, 4).NumberFormat = $ws.Cells.Item(516, 4).NumberFormat
